# ComunicELA-tfg - tableros/verbos.xlsx
# "UI cambiada y arreglados fallitos"
#
# The board had the verb "IR" in cell C4; replace it with "QUERER" and
# move the active selection to C4 (matching where the edit was made).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "QUERER"
$ws.Range("C4").Select()
